$d = $word.ActiveDocument

# Update the date in the title paragraph (unique text, safe to use Find/Replace).
$d.Content.Find.Execute("2025-11-05 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-06 Thursday", 2) | Out-Null

# Update each table cell by explicit (row, column) address, assigning the new
# text straight onto the cell Range. This avoids ambiguity from duplicate old
# values living in different cells (Find/Replace isn't reliably range-scoped).
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "69÷5=13, 4"
$tbl.Cell(1, 2).Range.Text = "72÷4=18, 0"
$tbl.Cell(1, 3).Range.Text = "66÷8=8, 2"
$tbl.Cell(1, 4).Range.Text = "27÷7=3, 6"
$tbl.Cell(1, 5).Range.Text = "91÷5=18, 1"
$tbl.Cell(5, 1).Range.Text = "23÷5=4, 3"
$tbl.Cell(5, 2).Range.Text = "90÷6=15, 0"
$tbl.Cell(5, 3).Range.Text = "25÷8=3, 1"
$tbl.Cell(5, 4).Range.Text = "20÷8=2, 4"
$tbl.Cell(5, 5).Range.Text = "81÷6=13, 3"
$tbl.Cell(9, 1).Range.Text = "82÷7=11, 5"
$tbl.Cell(9, 2).Range.Text = "56÷7=8, 0"
$tbl.Cell(9, 3).Range.Text = "66÷3=22, 0"
$tbl.Cell(9, 4).Range.Text = "71÷5=14, 1"
$tbl.Cell(9, 5).Range.Text = "44÷9=4, 8"
$tbl.Cell(13, 1).Range.Text = "38÷7=5, 3"
$tbl.Cell(13, 2).Range.Text = "41÷5=8, 1"
$tbl.Cell(13, 3).Range.Text = "82÷8=10, 2"
$tbl.Cell(13, 4).Range.Text = "36÷9=4, 0"
$tbl.Cell(13, 5).Range.Text = "32÷2=16, 0"
$tbl.Cell(17, 1).Range.Text = "30÷7=4, 2"
$tbl.Cell(17, 2).Range.Text = "82÷7=11, 5"
$tbl.Cell(17, 3).Range.Text = "23÷6=3, 5"
$tbl.Cell(17, 4).Range.Text = "49÷7=7, 0"
$tbl.Cell(17, 5).Range.Text = "22÷6=3, 4"
